$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for all data rows 2-47 from 46070 to 46072
$ws.Range("C2:C47").Value = 46072

# Rows 18-47 got reshuffled (their A/B/G values rotated among the rows).
# Apply the new A (Beteckning), B (Datum), G (Area) values row by row,
# taken from the corresponding original row before the reshuffle.
$ws.Cells.Item(18, 1).Value = "A 22033-2025"
$ws.Cells.Item(18, 2).Value = 45785.29449074074
$ws.Cells.Item(18, 7).Value = 5.2
$ws.Cells.Item(19, 1).Value = "A 31321-2025"
$ws.Cells.Item(19, 2).Value = 45833
$ws.Cells.Item(19, 7).Value = 6.4
$ws.Cells.Item(20, 1).Value = "A 48265-2025"
$ws.Cells.Item(20, 2).Value = 45933
$ws.Cells.Item(20, 7).Value = 2.1
$ws.Cells.Item(21, 1).Value = "A 12977-2025"
$ws.Cells.Item(21, 2).Value = 45734.45465277778
$ws.Cells.Item(21, 7).Value = 2.1
$ws.Cells.Item(22, 1).Value = "A 61167-2024"
$ws.Cells.Item(22, 2).Value = 45645
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(23, 1).Value = "A 21536-2024"
$ws.Cells.Item(23, 2).Value = 45441.59925925926
$ws.Cells.Item(23, 7).Value = 2.8
$ws.Cells.Item(24, 1).Value = "A 44926-2025"
$ws.Cells.Item(24, 2).Value = 45918.55856481481
$ws.Cells.Item(24, 7).Value = 3.8
$ws.Cells.Item(25, 1).Value = "A 43067-2024"
$ws.Cells.Item(25, 2).Value = 45567.47446759259
$ws.Cells.Item(25, 7).Value = 1.1
$ws.Cells.Item(26, 1).Value = "A 12273-2024"
$ws.Cells.Item(26, 2).Value = 45378.47817129629
$ws.Cells.Item(26, 7).Value = 0.9
$ws.Cells.Item(27, 1).Value = "A 33201-2023"
$ws.Cells.Item(27, 2).Value = 45127.42379629629
$ws.Cells.Item(27, 7).Value = 0.9
$ws.Cells.Item(28, 1).Value = "A 6042-2024"
$ws.Cells.Item(28, 2).Value = 45336
$ws.Cells.Item(28, 7).Value = 1.7
$ws.Cells.Item(29, 1).Value = "A 37570-2025"
$ws.Cells.Item(29, 2).Value = 45880.37358796296
$ws.Cells.Item(29, 7).Value = 0.9
$ws.Cells.Item(30, 1).Value = "A 22953-2023"
$ws.Cells.Item(30, 2).Value = 45072
$ws.Cells.Item(30, 7).Value = 1.9
$ws.Cells.Item(31, 1).Value = "A 62804-2023"
$ws.Cells.Item(31, 2).Value = 45271
$ws.Cells.Item(31, 7).Value = 0.6
$ws.Cells.Item(32, 1).Value = "A 62768-2025"
$ws.Cells.Item(32, 2).Value = 46008.59856481481
$ws.Cells.Item(32, 7).Value = 4.2
$ws.Cells.Item(33, 1).Value = "A 3811-2024"
$ws.Cells.Item(33, 2).Value = 45321.673125
$ws.Cells.Item(33, 7).Value = 0.9
$ws.Cells.Item(34, 1).Value = "A 61178-2024"
$ws.Cells.Item(34, 2).Value = 45645
$ws.Cells.Item(34, 7).Value = 6.6
$ws.Cells.Item(35, 1).Value = "A 46384-2025"
$ws.Cells.Item(35, 2).Value = 45925
$ws.Cells.Item(35, 7).Value = 3.1
$ws.Cells.Item(36, 1).Value = "A 3676-2022"
$ws.Cells.Item(36, 2).Value = 44586
$ws.Cells.Item(36, 7).Value = 0.5
$ws.Cells.Item(37, 1).Value = "A 6258-2024"
$ws.Cells.Item(37, 2).Value = 45337.77947916667
$ws.Cells.Item(37, 7).Value = 1.4
$ws.Cells.Item(38, 1).Value = "A 46379-2025"
$ws.Cells.Item(38, 2).Value = 45925
$ws.Cells.Item(38, 7).Value = 7.1
$ws.Cells.Item(39, 1).Value = "A 88-2025"
$ws.Cells.Item(39, 2).Value = 45659.46386574074
$ws.Cells.Item(39, 7).Value = 1.4
$ws.Cells.Item(40, 1).Value = "A 8639-2023"
$ws.Cells.Item(40, 2).Value = 44977.95614583333
$ws.Cells.Item(40, 7).Value = 3.8
$ws.Cells.Item(41, 1).Value = "A 60809-2024"
$ws.Cells.Item(41, 2).Value = 45644.61414351852
$ws.Cells.Item(41, 7).Value = 0.5
$ws.Cells.Item(42, 1).Value = "A 11256-2024"
$ws.Cells.Item(42, 2).Value = 45371.66233796296
$ws.Cells.Item(42, 7).Value = 0.6
$ws.Cells.Item(43, 1).Value = "A 24771-2023"
$ws.Cells.Item(43, 2).Value = 45084.64277777778
$ws.Cells.Item(43, 7).Value = 1
$ws.Cells.Item(44, 1).Value = "A 56133-2023"
$ws.Cells.Item(44, 2).Value = 45240
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(45, 1).Value = "A 11261-2024"
$ws.Cells.Item(45, 2).Value = 45371.67425925926
$ws.Cells.Item(45, 7).Value = 2.7
$ws.Cells.Item(46, 1).Value = "A 19295-2025"
$ws.Cells.Item(46, 2).Value = 45769.56212962963
$ws.Cells.Item(46, 7).Value = 2.8
$ws.Cells.Item(47, 1).Value = "A 21477-2025"
$ws.Cells.Item(47, 2).Value = 45782.59390046296
$ws.Cells.Item(47, 7).Value = 2.9
